$d = $word.ActiveDocument

# Locate the paragraph that holds the "<<date_placeholder>>" token (it may be
# split across several runs, e.g. "<<" / "date" / "_placeholder" / ">>").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*date_placeholder*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs(1)
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Rebuild the paragraph as two Arial paragraphs:
#   "Date: <<date_placeholder>>"
#   "Location: <<location_placeholder>>"
# The first paragraph also carries the Arial font on the paragraph mark
# itself (w:pPr/w:rPr), matching a whole-paragraph font change in Word.
$newXml = '<w:p ' + $wNs + '>' +
            '<w:pPr>' +
              '<w:rPr>' +
                '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '</w:rPr>' +
            '</w:pPr>' +
            '<w:r>' +
              '<w:rPr>' +
                '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '</w:rPr>' +
              '<w:t>Date: &lt;&lt;date_placeholder&gt;&gt;</w:t>' +
            '</w:r>' +
          '</w:p>' +
          '<w:p ' + $wNs + '>' +
            '<w:r>' +
              '<w:rPr>' +
                '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
              '</w:rPr>' +
              '<w:t>Location: &lt;&lt;location_placeholder&gt;&gt;</w:t>' +
            '</w:r>' +
          '</w:p>'

$target.Range.InsertXML($newXml)
